$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new PANC / Prostate v1.2-consortium rows (2023-01 release) ---
$ws.Cells.Item(117,1).Value = "PANC"
$ws.Cells.Item(117,2).Value = "v1.2-consortium"
$ws.Cells.Item(117,3).Value = "ca_radtx_dataset"
$ws.Cells.Item(117,5).Value = "2023-01"
$ws.Cells.Item(117,4).Value = "syn50908650"
$ws.Cells.Item(118,1).Value = "PANC"
$ws.Cells.Item(118,2).Value = "v1.2-consortium"
$ws.Cells.Item(118,3).Value = "cancer_level_dataset_index"
$ws.Cells.Item(118,5).Value = "2023-01"
$ws.Cells.Item(118,4).Value = "syn50908651"
$ws.Cells.Item(119,1).Value = "PANC"
$ws.Cells.Item(119,2).Value = "v1.2-consortium"
$ws.Cells.Item(119,3).Value = "cancer_level_dataset_non_index"
$ws.Cells.Item(119,5).Value = "2023-01"
$ws.Cells.Item(119,4).Value = "syn50908652"
$ws.Cells.Item(120,1).Value = "PANC"
$ws.Cells.Item(120,2).Value = "v1.2-consortium"
$ws.Cells.Item(120,3).Value = "cancer_panel_test_level_dataset"
$ws.Cells.Item(120,5).Value = "2023-01"
$ws.Cells.Item(120,4).Value = "syn50908653"
$ws.Cells.Item(121,1).Value = "PANC"
$ws.Cells.Item(121,2).Value = "v1.2-consortium"
$ws.Cells.Item(121,3).Value = "imaging_level_dataset"
$ws.Cells.Item(121,5).Value = "2023-01"
$ws.Cells.Item(121,4).Value = "syn50908654"
$ws.Cells.Item(122,1).Value = "PANC"
$ws.Cells.Item(122,2).Value = "v1.2-consortium"
$ws.Cells.Item(122,3).Value = "med_onc_note_level_dataset"
$ws.Cells.Item(122,5).Value = "2023-01"
$ws.Cells.Item(122,4).Value = "syn50908655"
$ws.Cells.Item(123,1).Value = "PANC"
$ws.Cells.Item(123,2).Value = "v1.2-consortium"
$ws.Cells.Item(123,3).Value = "pathology_report_level_dataset"
$ws.Cells.Item(123,5).Value = "2023-01"
$ws.Cells.Item(123,4).Value = "syn50908656"
$ws.Cells.Item(124,1).Value = "PANC"
$ws.Cells.Item(124,2).Value = "v1.2-consortium"
$ws.Cells.Item(124,3).Value = "patient_level_dataset"
$ws.Cells.Item(124,5).Value = "2023-01"
$ws.Cells.Item(124,4).Value = "syn50908657"
$ws.Cells.Item(125,1).Value = "PANC"
$ws.Cells.Item(125,2).Value = "v1.2-consortium"
$ws.Cells.Item(125,3).Value = "regimen_cancer_level_dataset"
$ws.Cells.Item(125,5).Value = "2023-01"
$ws.Cells.Item(125,4).Value = "syn50908658"
$ws.Cells.Item(126,1).Value = "PANC"
$ws.Cells.Item(126,2).Value = "v1.2-consortium"
$ws.Cells.Item(126,3).Value = "tm_level_dataset"
$ws.Cells.Item(126,5).Value = "2023-01"
$ws.Cells.Item(126,4).Value = "syn50908659"
$ws.Cells.Item(127,1).Value = "Prostate"
$ws.Cells.Item(127,2).Value = "v1.2-consortium"
$ws.Cells.Item(127,3).Value = "ca_radtx_dataset"
$ws.Cells.Item(127,5).Value = "2023-01"
$ws.Cells.Item(127,4).Value = "syn50908660"
$ws.Cells.Item(128,1).Value = "Prostate"
$ws.Cells.Item(128,2).Value = "v1.2-consortium"
$ws.Cells.Item(128,3).Value = "cancer_level_dataset_index"
$ws.Cells.Item(128,5).Value = "2023-01"
$ws.Cells.Item(128,4).Value = "syn50908661"
$ws.Cells.Item(129,1).Value = "Prostate"
$ws.Cells.Item(129,2).Value = "v1.2-consortium"
$ws.Cells.Item(129,3).Value = "cancer_level_dataset_non_index"
$ws.Cells.Item(129,5).Value = "2023-01"
$ws.Cells.Item(129,4).Value = "syn50908662"
$ws.Cells.Item(130,1).Value = "Prostate"
$ws.Cells.Item(130,2).Value = "v1.2-consortium"
$ws.Cells.Item(130,3).Value = "cancer_panel_test_level_dataset"
$ws.Cells.Item(130,5).Value = "2023-01"
$ws.Cells.Item(130,4).Value = "syn50908663"
$ws.Cells.Item(131,1).Value = "Prostate"
$ws.Cells.Item(131,2).Value = "v1.2-consortium"
$ws.Cells.Item(131,3).Value = "imaging_level_dataset"
$ws.Cells.Item(131,5).Value = "2023-01"
$ws.Cells.Item(131,4).Value = "syn50908664"
$ws.Cells.Item(132,1).Value = "Prostate"
$ws.Cells.Item(132,2).Value = "v1.2-consortium"
$ws.Cells.Item(132,3).Value = "med_onc_note_level_dataset"
$ws.Cells.Item(132,5).Value = "2023-01"
$ws.Cells.Item(132,4).Value = "syn50908665"
$ws.Cells.Item(133,1).Value = "Prostate"
$ws.Cells.Item(133,2).Value = "v1.2-consortium"
$ws.Cells.Item(133,3).Value = "pathology_report_level_dataset"
$ws.Cells.Item(133,5).Value = "2023-01"
$ws.Cells.Item(133,4).Value = "syn50908666"
$ws.Cells.Item(134,1).Value = "Prostate"
$ws.Cells.Item(134,2).Value = "v1.2-consortium"
$ws.Cells.Item(134,3).Value = "patient_level_dataset"
$ws.Cells.Item(134,5).Value = "2023-01"
$ws.Cells.Item(134,4).Value = "syn50908667"
$ws.Cells.Item(135,1).Value = "Prostate"
$ws.Cells.Item(135,2).Value = "v1.2-consortium"
$ws.Cells.Item(135,3).Value = "regimen_cancer_level_dataset"
$ws.Cells.Item(135,5).Value = "2023-01"
$ws.Cells.Item(135,4).Value = "syn50908668"
$ws.Cells.Item(136,1).Value = "Prostate"
$ws.Cells.Item(136,2).Value = "v1.2-consortium"
$ws.Cells.Item(136,3).Value = "tm_level_dataset"
$ws.Cells.Item(136,5).Value = "2023-01"
$ws.Cells.Item(136,4).Value = "syn50908669"

# --- Update existing cell D105: regimen_cancer_level_dataset synapse id for BrCa v1.2-consortium ---
$ws.Cells.Item(105,1).Value = "BrCa"
$ws.Cells.Item(105,2).Value = "v1.2-consortium"
$ws.Cells.Item(105,3).Value = "regimen_cancer_level_dataset"
$ws.Cells.Item(105,4).Value = "syn43172837"
$ws.Cells.Item(105,5).Value = "2022-10"

# --- Update sheet view: remove topLeftCell scroll position, change selection ---
$ws.Range("K21").Select() | Out-Null

# --- Update workbook window geometry to match target (best effort) ---
$win = $excel.ActiveWindow
$win.Left = -120
$win.Top = -120
$win.Width = 29040
$win.Height = 17640
